$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row above row 2 so all existing records (rows 2-18) shift down
# to rows 3-19, making room for the new weekly price observation.
$ws.Rows.Item(2).Insert()

# The insert operation can inherit formatting from the header row above;
# reset the new row's formatting to match the other (now-shifted) data
# rows by copying the format from row 3 (the former row 2).
$ws.Range("A3:R3").Copy()
$ws.Range("A2:R2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new record's fixed columns (same market/product as all others).
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(2, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(2, 4).Value = 45257
$ws.Cells.Item(2, 5).Value = 15
$ws.Cells.Item(2, 6).Value = 100112030
$ws.Cells.Item(2, 7).Value = "Poroto granado"
$ws.Cells.Item(2, 8).Value = "Sin especificar"
$ws.Cells.Item(2, 9).Value = "Primera"
$ws.Cells.Item(2, 10).Value = 1000
$ws.Cells.Item(2, 11).Value = 1100
$ws.Cells.Item(2, 12).Value = 1200
$ws.Cells.Item(2, 13).Value = 1150
$ws.Cells.Item(2, 14).Value = "`$/kilo"
$ws.Cells.Item(2, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(2, 16).Value = 1150
$ws.Cells.Item(2, 17).Value = 1
$ws.Cells.Item(2, 18).Value = "Hortaliza"
